$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.403.22'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').Value = '3.686.17'
$ws.Range('E3').Value = '  -2.96%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '682.12'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.60'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.16%  '
$ws.Range('D7').Value = '3.685.51'
$ws.Range('E7').Value = '  -2.96%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -4.14%  '
$ws.Range('E10').Value = '  -7.60%  '
$ws.Range('E11').Value = '  -3.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('E13').Value = '  -4.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.52'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -6.05%  '
$ws.Range('D15').Value = '4.308.61'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('D16').Value = '3.686.93'
$ws.Range('E16').Value = '  -3.33%  '
$ws.Range('D17').Value = '69.425.18'
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('E19').Value = '  -6.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.63'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '482.12'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.91'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -6.99%  '
$ws.Range('E23').Value = '  -7.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.36'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.49%  '
$ws.Range('D25').Value = '3.832.04'
$ws.Range('E25').Value = '  -2.94%  '
$ws.Range('E26').Value = '  -8.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.50'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.22%  '
$ws.Range('E29').Value = '  -6.45%  '
$ws.Range('E30').Value = '  -8.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.72'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -9.42%  '
$ws.Range('E32').Value = '  -7.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.86'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -6.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.10'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.39%  '
$ws.Range('E35').Value = '  -4.65%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '3.653.78'
$ws.Range('E37').Value = '  -3.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.52'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -5.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.37'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0937'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -7.28%  '
$ws.Range('E41').Value = '  -4.65%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.956'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '159.64'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '48.39'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('E47').Value = '  -11.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.04'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +6.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000289'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -8.27%  '
$ws.Range('E50').Value = '  +1.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '394.58'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.20%  '
